# Updated the bay numbering system to maintain compatibility with the existing data structure

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Flightline 3 -> 1, Serial Number "teet" -> "56382", Customer Name "t4etw" -> "Hass"
$ws.Range("B2").Value = 1
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "56382"
$ws.Range("D2").Value = "Hass"

# Row 3: Bay Number 2 -> 3, Flightline 2 -> 1, Serial Number "Test" -> "56887",
# Customer Name "Test" -> "hass bombn"
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "56887"
$ws.Range("D3").Value = "hass bombn"

# The old rows 4-8 (Bell IT, Test2dwadadw, test2, fssd, rwef) are no longer
# part of the bay list, so remove them and let the used range shrink to A1:F3
$ws.Range("4:8").EntireRow.Delete()

# Give the columns their explicit widths (Excel's ColumnWidth property is
# offset by 5/6 of a character from the stored <col width> attribute, so
# subtract that to land on the target widths of 10, 10, 15, 20, 10, 50)
$ws.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 49.166666666666664
